$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 0. The source workbook contains a number of cells that are declared but
#    carry no actual value (empty inline strings). Clear them out so they
#    are dropped entirely when the sheet is re-serialized.
$emptyCells = @("S2", "AA2", "AN2", "AP2", "AR2", "AS2", "N3", "Q3", "R3", "S3", "T3", "U3", "V3", "W3", "X3", "Y3", "Z3", "AA3", "AN3", "AO3", "AP3", "AQ3", "AR3", "AS3", "S4", "AA4", "AN4", "AP4", "AR4", "AS4", "S5", "U5", "V5", "AA5", "AN5", "AP5", "AR5", "AS5", "S6", "AN6", "AR6", "AS6", "S7", "AN7", "AQ7", "AR7", "AS7", "S8", "AA8", "AN8", "AQ8", "AR8", "AS8", "AN9", "AQ9", "AR9", "AS9", "AN10", "AQ10", "AR10", "AS10", "AN11", "AO11", "AQ11", "AR11")
foreach ($addr in $emptyCells) {
    $ws.Range($addr).ClearContents()
}

# 1. Add the new "DropdownOptions" helper sheet right after Sheet1.
$dropdownSheet = $wb.Worksheets.Add($null, $ws)
$dropdownSheet.Name = "DropdownOptions"

$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $cell = $dropdownSheet.Cells.Item($i + 1, 1)
    # Prefix with an apostrophe so values like "100%" are stored as text
    # instead of being interpreted/converted into a percentage number.
    $cell.Value = "'" + $options[$i]
    $cell.Style = "Normal"
}

# Hide the helper sheet.
$dropdownSheet.Visible = $false

# 2. Add the new "Status as of July 4, 2025" column (AU) to Sheet1.
$ws.Range("AU1").Value = "Status as of July 4, 2025"

# 3. Apply a dropdown-list data validation on AU2:AU11 sourced from the
#    DropdownOptions helper sheet.
$validation = $ws.Range("AU2:AU11").Validation
$validation.Add(3, 1, 1, '=DropdownOptions!$A$1:$A$7')
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $false
$validation.ShowError = $false

$ws.Activate()
